$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F10").Value = 6901
$ws1.Range("F12").Value = 385
$ws1.Range("F13").Value = 3259

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F13").Value = 6901
$ws4.Range("F16").Value = 385
$ws4.Range("F17").Value = 3259
